# Actualización automática de tasas-transfi.xlsx

$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet Hoja1 (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.02 = 24855.54 pesos`n✅ 24855.54 pesos = 6.02 = 967.21 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $text

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 166
$wsTasas.Range("O10").Value = 4126.02
$wsTasas.Range("N12").Value = 4130
$wsTasas.Range("O12").Value = 160.712
